$d = $word.ActiveDocument

# --- Change 1: remove the empty paragraph that follows "Changes 2023-08-16" ---
# Find the paragraph containing the heading text, then drop the very next
# (empty) paragraph, paragraph mark included, so the following bullet list
# item becomes adjacent to the heading.
$changesPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd("`r", "`a") -eq "Changes 2023-08-16") {
        $changesPara = $p
        break
    }
}
if ($changesPara -ne $null) {
    $emptyPara = $changesPara.Next()
    if ($emptyPara -ne $null -and $emptyPara.Range.Text.Trim("`r", "`a") -eq "") {
        $emptyPara.Range.Delete()
    }
}

# --- Change 2: suppress automatic hyphenation on the built-in Normal style ---
$normalStyle = $d.Styles("Normal")
$normalStyle.ParagraphFormat.Hyphenation = $false
